# Paraguay Division Profesional - data refresh
# - Swap a handful of mis-ordered match rows (the underlying match ids were
#   transposed between consecutive fixtures) so each row again holds its own
#   match's data.
# - Append six new fixtures (ids 239-244) that were added to the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($r1, $r2)
    $rng1 = $ws.Range("B$r1" + ":AB$r1")
    $rng2 = $ws.Range("B$r2" + ":AB$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Pairs of rows whose match data (everything except the running "id" in
# column A) needs to be swapped back into the correct row.
Swap-Rows 3 4
Swap-Rows 130 131
Swap-Rows 140 141
Swap-Rows 143 145
Swap-Rows 236 237

# Append the six new fixtures as rows 241-246, reusing the formatting of the
# last existing data row (240).
$srcRange = $ws.Range("A240:AB240")
for ($r = 241; $r -le 246; $r++) {
    $destRange = $ws.Range("A$r" + ":AB$r")
    $srcRange.Copy($destRange)
}

$newRows = @(
    @{ Row=241; A=239; B=7609666; D=45422.79166666666; E="2 de Mayo";             F="Tacuary";               G=3; H=1; I="H";
       J=1.833; K=3.5;  L=4.333; M=1.85;  N=3.3;   O=4.333; P=-0.5;  Q=1.9;   R=1.9;   S=2.25; T=1.95;  U=1.85;
       V=0.8500000000000001; W=-1; X=-1; Y=0.8999999999999999; Z=-1; AA=0.95; AB=-1 },
    @{ Row=242; A=240; B=7609211; D=45422.89583333334; E="Libertad Asuncion";     F="Sol de America";        G=4; H=1; I="H";
       J=1.5;   K=4;    L=6.5;   M=1.444; N=4.2;   O=7.5;   P=-1.25; Q=2.025; R=1.775; S=2.5;  T=1.85;  U=1.95;
       V=0.444; W=-1; X=-1; Y=1.025; Z=-1; AA=0.8500000000000001; AB=-1 },
    @{ Row=243; A=241; B=7609168; D=45423.75;          E="Sportivo Trinidense";  F="Nacional Asuncion";     G=0; H=1; I="A";
       J=2.7;   K=3.2;  L=2.55;  M=2.55;  N=3.2;   O=2.7;   P=0;     Q=1.825; R=1.975; S=2.25; T=1.775; U=2.025;
       V=-1; W=-1; X=1.7; Y=-1; Z=0.9750000000000001; AA=-1; AB=1.025 },
    @{ Row=244; A=242; B=7609210; D=45423.85416666666; E="Sportivo Luqueno";     F="Sportivo Ameliano";     G=0; H=0; I="D";
       J=1.833; K=3.5;  L=4.2;   M=1.6;   N=3.75;  O=6;     P=-0.75; Q=1.775; R=2.025; S=2.25; T=1.8;   U=2;
       V=-1; W=2.75; X=-1; Y=-1; Z=1.025; AA=-1; AB=1 },
    @{ Row=245; A=243; B=7609166; D=45424.72916666666; E="Cerro Porteno";        F="Olimpia Asuncion";      G=1; H=1; I="D";
       J=1.909; K=3.4;  L=4;     M=1.7;   N=3.6;   O=5.25;  P=-0.75; Q=1.9;   R=1.9;   S=2.25; T=1.925; U=1.875;
       V=-1; W=2.6; X=-1; Y=-1; Z=0.8999999999999999; AA=-0.5; AB=0.4375 },
    @{ Row=246; A=244; B=7609167; D=45425.85416666666; E="Guarani Asuncion";     F="General Caballero JLM"; G=0; H=0; I="D";
       J=1.7;   K=4;    L=4.5;   M=1.65;  N=4.333; O=4.5;   P=-0.75; Q=1.8;   R=2;     S=2.75; T=1.775; U=2.025;
       V=-1; W=3.333; X=-1; Y=-1; Z=1; AA=-1; AB=1.025 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Range("A$r").Value2 = $row.A
    $ws.Range("B$r").Value2 = $row.B
    $ws.Range("C$r").Value = "Paraguay Division Profesional"
    $ws.Range("D$r").Value2 = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value2 = $row.G
    $ws.Range("H$r").Value2 = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value2 = $row.J
    $ws.Range("K$r").Value2 = $row.K
    $ws.Range("L$r").Value2 = $row.L
    $ws.Range("M$r").Value2 = $row.M
    $ws.Range("N$r").Value2 = $row.N
    $ws.Range("O$r").Value2 = $row.O
    $ws.Range("P$r").Value2 = $row.P
    $ws.Range("Q$r").Value2 = $row.Q
    $ws.Range("R$r").Value2 = $row.R
    $ws.Range("S$r").Value2 = $row.S
    $ws.Range("T$r").Value2 = $row.T
    $ws.Range("U$r").Value2 = $row.U
    $ws.Range("V$r").Value2 = $row.V
    $ws.Range("W$r").Value2 = $row.W
    $ws.Range("X$r").Value2 = $row.X
    $ws.Range("Y$r").Value2 = $row.Y
    $ws.Range("Z$r").Value2 = $row.Z
    $ws.Range("AA$r").Value2 = $row.AA
    $ws.Range("AB$r").Value2 = $row.AB
}

Write-Output "Applied Paraguay Division Profesional update"
